# Auto-generated edit script applying the Mateus_Profits.xlsx data-refresh diff
# Sheets (workbook tabs) map to the diff's per-file hunks: ALC, ARM, BSM, CRP, CUL, LTW, WVR
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1351.5
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H8").Value = 33
$ws.Range("I8").Value = 33
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 99
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 40
$ws.Range("N8").ClearContents()
$ws.Range("H9").Value = 169.24324
$ws.Range("J9").Value = 330
$ws.Range("L9").Value = 330
$ws.Range("N9").Value = -668
$ws.Range("H32").Value = 6045.6
$ws.Range("J32").Value = 6974.25
$ws.Range("L32").Value = 6974.25
$ws.Range("N32").Value = -7626.25
$ws.Range("H106").Value = 7265.6665
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 668.5714
$ws.Range("I107").Value = 785.8182
$ws.Range("J107").Value = 238.66667
$ws.Range("K107").Value = 785.8182
$ws.Range("L107").Value = 238.66667
$ws.Range("M107").Value = 1134.1818
$ws.Range("N107").Value = -4078.66667
$ws.Range("H112").Value = 2284.889
$ws.Range("J112").Value = 2112.7144
$ws.Range("L112").Value = 6338.1432
$ws.Range("N112").Value = -8554.143199999999
$ws.Range("H132").Value = 2094.842
$ws.Range("I132").Value = 1668.5625
$ws.Range("K132").Value = 5005.6875
$ws.Range("M132").Value = -2475.6875
$ws.Range("H135").Value = 2971
$ws.Range("I135").Value = 2091.2856
$ws.Range("K135").Value = 18821.5704
$ws.Range("M135").Value = -16286.5704
$ws.Range("H138").Value = 3788.973
$ws.Range("J138").Value = 3948.0476
$ws.Range("L138").Value = 11844.1428
$ws.Range("N138").Value = -22124.1428

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4202.079
$ws.Range("I32").Value = 4149.4507
$ws.Range("J32").Value = 4949.4
$ws.Range("K32").Value = 4149.4507
$ws.Range("L32").Value = 4949.4
$ws.Range("M32").Value = -3862.4507
$ws.Range("N32").Value = -5523.4
$ws.Range("H102").Value = 3875.2856
$ws.Range("I102").Value = 2521.1667
$ws.Range("K102").Value = 2521.1667
$ws.Range("M102").Value = -899.1667000000002
$ws.Range("H111").Value = 644
$ws.Range("J111").Value = 644
$ws.Range("L111").Value = 644
$ws.Range("N111").Value = -8824
$ws.Range("H122").Value = 3766.8572
$ws.Range("I122").Value = 3037.8635
$ws.Range("J122").Value = 6439.8335
$ws.Range("K122").Value = 9113.5905
$ws.Range("L122").Value = 19319.5005
$ws.Range("M122").Value = -6663.5905
$ws.Range("N122").Value = -24219.5005
$ws.Range("H132").Value = 4188.317
$ws.Range("I132").Value = 3797.2354
$ws.Range("K132").Value = 11391.7062
$ws.Range("M132").Value = -8861.706200000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H105").Value = 3959.7693
$ws.Range("I105").Value = 3548.3
$ws.Range("J105").Value = 5331.3335
$ws.Range("K105").Value = 3548.3
$ws.Range("L105").Value = 5331.3335
$ws.Range("M105").Value = -1801.3
$ws.Range("N105").Value = -8825.333500000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2619.8076
$ws.Range("I16").Value = 2150.7273
$ws.Range("K16").Value = 2150.7273
$ws.Range("M16").Value = -1863.7273
$ws.Range("H31").Value = 4434.028
$ws.Range("I31").Value = 3380.7083
$ws.Range("K31").Value = 3380.7083
$ws.Range("M31").Value = -3085.7083
$ws.Range("H34").Value = 4434.028
$ws.Range("I34").Value = 3380.7083
$ws.Range("K34").Value = 3380.7083
$ws.Range("M34").Value = -3178.7083
$ws.Range("H41").Value = 17950
$ws.Range("I41").Value = 19800
$ws.Range("J41").Value = 16562.5
$ws.Range("K41").Value = 19800
$ws.Range("L41").Value = 16562.5
$ws.Range("M41").Value = -19372
$ws.Range("N41").Value = -17418.5
$ws.Range("H86").Value = 10492.5
$ws.Range("I86").Value = 5985
$ws.Range("J86").Value = 15000
$ws.Range("K86").Value = 5985
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = -4862
$ws.Range("N86").Value = -17246
$ws.Range("H89").Value = 10492.5
$ws.Range("I89").Value = 5985
$ws.Range("J89").Value = 15000
$ws.Range("K89").Value = 29925
$ws.Range("L89").Value = 75000
$ws.Range("M89").Value = -24309
$ws.Range("N89").Value = -86232
$ws.Range("H113").Value = 2619.8076
$ws.Range("I113").Value = 2150.7273
$ws.Range("K113").Value = 2150.7273
$ws.Range("M113").Value = 19.27269999999999
$ws.Range("H132").Value = 3917.5
$ws.Range("I132").Value = 3063.4
$ws.Range("J132").Value = 6479.8
$ws.Range("K132").Value = 9190.200000000001
$ws.Range("L132").Value = 19439.4
$ws.Range("M132").Value = -6660.200000000001
$ws.Range("N132").Value = -24499.4
$ws.Range("H141").Value = 39333.766
$ws.Range("J141").Value = 40042.125
$ws.Range("L141").Value = 40042.125
$ws.Range("N141").Value = -50402.125

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 183718.64
$ws.Range("I26").Value = 285886.44
$ws.Range("K26").Value = 857659.3200000001
$ws.Range("M26").Value = -857371.3200000001
$ws.Range("H64").Value = 4999
$ws.Range("I64").Value = 4999
$ws.Range("K64").Value = 14997
$ws.Range("M64").Value = -14727
$ws.Range("H67").Value = 4999
$ws.Range("I67").Value = 4999
$ws.Range("K67").Value = 14997
$ws.Range("M67").Value = -14061
$ws.Range("H81").Value = 10507
$ws.Range("I81").Value = 2260.5
$ws.Range("K81").Value = 6781.5
$ws.Range("M81").Value = -5658.5
$ws.Range("H84").Value = 10507
$ws.Range("I84").Value = 2260.5
$ws.Range("K84").Value = 20344.5
$ws.Range("M84").Value = -14728.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15560.211
$ws.Range("I7").Value = 11711.8
$ws.Range("K7").Value = 11711.8
$ws.Range("M7").Value = -11599.8
$ws.Range("H46").Value = 7685.2354
$ws.Range("J46").Value = 11457.417
$ws.Range("L46").Value = 11457.417
$ws.Range("N46").Value = -11833.417
$ws.Range("H55").Value = 813.6667
$ws.Range("I55").Value = 983.5
$ws.Range("K55").Value = 983.5
$ws.Range("M55").Value = -810.5
$ws.Range("H61").Value = 47962.043
$ws.Range("I61").Value = 54156.2
$ws.Range("J61").Value = 6667.6665
$ws.Range("K61").Value = 54156.2
$ws.Range("L61").Value = 6667.6665
$ws.Range("M61").Value = -53954.2
$ws.Range("N61").Value = -7071.6665
$ws.Range("H113").Value = 47962.043
$ws.Range("I113").Value = 54156.2
$ws.Range("J113").Value = 6667.6665
$ws.Range("K113").Value = 54156.2
$ws.Range("L113").Value = 6667.6665
$ws.Range("M113").Value = -51986.2
$ws.Range("N113").Value = -11007.6665
$ws.Range("H126").Value = 15560.211
$ws.Range("I126").Value = 11711.8
$ws.Range("K126").Value = 35135.39999999999
$ws.Range("M126").Value = -32665.39999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1900
$ws.Range("I100").Value = 1900
$ws.Range("K100").Value = 3800
$ws.Range("M100").Value = -3259
$ws.Range("H107").Value = 981.5
$ws.Range("J107").Value = 750
$ws.Range("L107").Value = 2250
$ws.Range("N107").Value = -6090
$ws.Range("H113").Value = 488.2857
$ws.Range("I113").Value = 490.47827
$ws.Range("J113").Value = 478.2
$ws.Range("K113").Value = 1471.43481
$ws.Range("L113").Value = 1434.6
$ws.Range("M113").Value = 698.56519
$ws.Range("N113").Value = -5774.6
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 3097.95
$ws.Range("I122").Value = 2119.7144
$ws.Range("K122").Value = 6359.1432
$ws.Range("M122").Value = -3909.1432
$ws.Range("H132").Value = 4069.2903
$ws.Range("I132").Value = 2949.963
$ws.Range("J132").Value = 11624.75
$ws.Range("K132").Value = 8849.889000000001
$ws.Range("L132").Value = 34874.25
$ws.Range("M132").Value = -6319.889000000001
$ws.Range("N132").Value = -39934.25
$ws.Range("H141").Value = 65358
$ws.Range("J141").Value = 65358
$ws.Range("L141").Value = 65358
$ws.Range("N141").Value = -75718
